# Add 96 new "package" rows (pack1..pack96) to the "packages" worksheet.
# These land in column A starting at row 5 (rows 1-4 already hold the
# header + the existing "test" / "molgenis" / "org" package rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("packages")

for ($i = 1; $i -le 96; $i++) {
    $row = 4 + $i
    $ws.Cells.Item($row, 1).Value = "pack$i"
}

# Match the author's final selection on that sheet.
$ws.Activate()
$ws.Range("A2").Select()
